$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
  $cell = $ws.Range($addr)
  $origStyle = $cell.Style
  $cell.NumberFormat = "@"
  $cell.Value = $value
  $cell.Style = $origStyle
}

# Row 2
Set-TextValue 'D2' '30.685.87'
Set-TextValue 'E2' '  +1.65%  '

# Row 3
Set-TextValue 'D3' '1.893.34'
Set-TextValue 'E3' '  +2.16%  '

# Row 4
Set-TextValue 'D4' '0.9998'
Set-TextValue 'E4' '  -0.13%  '

# Row 5
Set-TextValue 'D5' '238.63'
Set-TextValue 'E5' '  +1.40%  '

# Row 6
Set-TextValue 'D6' '0.9996'
Set-TextValue 'E6' '  -0.12%  '

# Row 7
Set-TextValue 'D7' '0.4843'
Set-TextValue 'E7' '  +1.47%  '

# Row 8
Set-TextValue 'D8' '0.2879'
Set-TextValue 'E8' '  +2.51%  '

# Row 9
Set-TextValue 'E9' '  +1.75%  '

# Row 10
Set-TextValue 'B10' 'Solana'
Set-TextValue 'C10' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D10' '16.88'
Set-TextValue 'E10' '  +4.27%  '

# Row 11
Set-TextValue 'B11' 'WrappedEther'
Set-TextValue 'C11' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D11' '1.831.43'
Set-TextValue 'E11' '  -1.13%  '

# Row 12
Set-TextValue 'D12' '0.07463'
Set-TextValue 'E12' '  +1.44%  '

# Row 13
Set-TextValue 'D13' '5.114'
Set-TextValue 'E13' '  +0.64%  '

# Row 14
Set-TextValue 'D14' '88.15'
Set-TextValue 'E14' '  +1.33%  '

# Row 15
Set-TextValue 'D15' '0.6696'
Set-TextValue 'E15' '  +3.89%  '

# Row 16
Set-TextValue 'D16' '30.666.58'
Set-TextValue 'E16' '  +1.77%  '

# Row 17
Set-TextValue 'D17' '13.27'
Set-TextValue 'E17' '  +1.20%  '

# Row 18
Set-TextValue 'D18' '1.0000'
Set-TextValue 'E18' '  -0.01%  '

# Row 19
Set-TextValue 'B19' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C19' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D19' '2.193.08'
Set-TextValue 'E19' '  +4.76%  '

# Row 20
Set-TextValue 'B20' 'ShibaInu'
Set-TextValue 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D20' '0.000007581'
Set-TextValue 'E20' '  +0.61%  '

# Row 21
Set-TextValue 'D21' '230.88'
Set-TextValue 'E21' '  +2.81%  '

# Row 22
Set-TextValue 'D22' '1.0000'
Set-TextValue 'E22' '  -0.14%  '

# Row 23
Set-TextValue 'E23' '  +0.18%  '

# Row 24
Set-TextValue 'D24' '6.212'
Set-TextValue 'E24' '  +2.39%  '

# Row 25
Set-TextValue 'D25' '170.10'
Set-TextValue 'E25' '  +4.45%  '

# Row 26
Set-TextValue 'D26' '9.387'
Set-TextValue 'E26' '  +2.08%  '

# Row 27
Set-TextValue 'E27' '  +2.18%  '

# Row 28
Set-TextValue 'D28' '1.964'
Set-TextValue 'E28' '  +2.60%  '

# Row 29
Set-TextValue 'D29' '0.1024'
Set-TextValue 'E29' '  +11.65%  '

# Row 30
Set-TextValue 'D30' '1.395'
Set-TextValue 'E30' '  -3.06%  '

# Row 31
Set-TextValue 'D31' '4.345'
Set-TextValue 'E31' '  +2.87%  '

# Row 32
Set-TextValue 'D32' '4.041'
Set-TextValue 'E32' '  +2.54%  '

# Row 33
Set-TextValue 'D33' '0.05064'
Set-TextValue 'E33' '  +2.12%  '

# Row 34
Set-TextValue 'E34' '  +6.88%  '

# Row 35
Set-TextValue 'D35' '0.7544'
Set-TextValue 'E35' '  +4.14%  '

# Row 36
Set-TextValue 'D36' '1.001'
Set-TextValue 'E36' '  +0.11%  '

# Row 37
Set-TextValue 'D37' '2.712'
Set-TextValue 'E37' '  +0.81%  '

# Row 38
Set-TextValue 'D38' '0.01881'
Set-TextValue 'E38' '  +2.69%  '

# Row 40
Set-TextValue 'D40' '0.9199'
Set-TextValue 'E40' '  +2.55%  '

# Row 41
Set-TextValue 'D41' '2.071'
Set-TextValue 'E41' '  +2.26%  '

# Row 42
Set-TextValue 'D42' '107.19'
Set-TextValue 'E42' '  +1.40%  '

# Row 43
Set-TextValue 'D43' '0.4305'
Set-TextValue 'E43' '  +1.92%  '

# Row 44
Set-TextValue 'E44' '  +0.28%  '

# Row 45
Set-TextValue 'D45' '5.682'
Set-TextValue 'E45' '  -4.17%  '

# Row 46
Set-TextValue 'D46' '7.435'
Set-TextValue 'E46' '  +1.34%  '

# Row 47
Set-TextValue 'D47' '64.66'
Set-TextValue 'E47' '  +1.31%  '

# Row 48
Set-TextValue 'D48' '0.1278'
Set-TextValue 'E48' '  -2.34%  '

# Row 49
Set-TextValue 'E49' '  -0.80%  '

# Row 50
Set-TextValue 'D50' '8.959'
Set-TextValue 'E50' '  +3.22%  '

# Row 51
Set-TextValue 'D51' '34.06'
Set-TextValue 'E51' '  +1.07%  '

